$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1374.5
$ws.Range("I2").Value = 1374.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1374.5
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1261.5
$ws.Range("N2").ClearContents()
$ws.Range("H9").Value = 244
$ws.Range("I9").Value = 214
$ws.Range("J9").Value = 274
$ws.Range("K9").Value = 214
$ws.Range("L9").Value = 274
$ws.Range("M9").Value = -45
$ws.Range("N9").Value = -612
$ws.Range("H43").Value = 6871.75
$ws.Range("I43").Value = 6662.3335
$ws.Range("K43").Value = 6662.3335
$ws.Range("M43").Value = -6593.3335
$ws.Range("H55").Value = 149.25
$ws.Range("I55").Value = 149.25
$ws.Range("K55").Value = 149.25
$ws.Range("M55").Value = 64.75
$ws.Range("H64").Value = 4000
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752
$ws.Range("H67").Value = 4000
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 15000
$ws.Range("M77").ClearContents()
$ws.Range("H98").Value = 664.8182
$ws.Range("I98").Value = 590.3333
$ws.Range("K98").Value = 590.3333
$ws.Range("M98").Value = 907.6667
$ws.Range("H106").Value = 8000
$ws.Range("J106").Value = 8000
$ws.Range("L106").Value = 8000
$ws.Range("N106").Value = -9262
$ws.Range("H112").Value = 1749.5555
$ws.Range("I112").Value = 850
$ws.Range("J112").Value = 1862
$ws.Range("K112").Value = 2550
$ws.Range("L112").Value = 5586
$ws.Range("M112").Value = -1442
$ws.Range("N112").Value = -7802
$ws.Range("H116").Value = 9147.727999999999
$ws.Range("I116").Value = 8788.6
$ws.Range("J116").Value = 9447
$ws.Range("K116").Value = 8788.6
$ws.Range("L116").Value = 9447
$ws.Range("M116").Value = -5346.6
$ws.Range("N116").Value = -16331
$ws.Range("H122").Value = 664.8182
$ws.Range("I122").Value = 590.3333
$ws.Range("K122").Value = 1770.9999
$ws.Range("M122").Value = 679.0001
$ws.Range("H137").Value = 2025.8334
$ws.Range("I137").Value = 1555.5
$ws.Range("K137").Value = 4666.5
$ws.Range("M137").Value = -2116.5
$ws.Range("H138").Value = 3315.386
$ws.Range("J138").Value = 3361.8164
$ws.Range("L138").Value = 10085.4492
$ws.Range("N138").Value = -20365.4492

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1054.9565
$ws.Range("I2").Value = 990.86664
$ws.Range("J2").Value = 1175.125
$ws.Range("K2").Value = 990.86664
$ws.Range("L2").Value = 1175.125
$ws.Range("M2").Value = -877.86664
$ws.Range("N2").Value = -1401.125
$ws.Range("H116").Value = 1054.9565
$ws.Range("I116").Value = 990.86664
$ws.Range("J116").Value = 1175.125
$ws.Range("K116").Value = 990.86664
$ws.Range("L116").Value = 1175.125
$ws.Range("M116").Value = 1303.13336
$ws.Range("N116").Value = -5763.125
$ws.Range("H132").Value = 2859.15
$ws.Range("I132").Value = 1816
$ws.Range("J132").Value = 4423.875
$ws.Range("K132").Value = 5448
$ws.Range("L132").Value = 13271.625
$ws.Range("M132").Value = -2918
$ws.Range("N132").Value = -18331.625
$ws.Range("H139").Value = 99999.5
$ws.Range("J139").Value = 99999.5
$ws.Range("L139").Value = 99999.5
$ws.Range("N139").Value = -110279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1054.9565
$ws.Range("I3").Value = 990.86664
$ws.Range("J3").Value = 1175.125
$ws.Range("K3").Value = 990.86664
$ws.Range("L3").Value = 1175.125
$ws.Range("M3").Value = -876.86664
$ws.Range("N3").Value = -1403.125
$ws.Range("H105").Value = 5826
$ws.Range("I105").Value = 4989
$ws.Range("J105").Value = 7500
$ws.Range("K105").Value = 4989
$ws.Range("L105").Value = 7500
$ws.Range("M105").Value = -3242
$ws.Range("N105").Value = -10994
$ws.Range("H134").Value = 2522.5386
$ws.Range("I134").Value = 2199.8572
$ws.Range("K134").Value = 6599.571599999999
$ws.Range("M134").Value = -4064.571599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2578.75
$ws.Range("I58").Value = 2083.2856
$ws.Range("K58").Value = 2083.2856
$ws.Range("M58").Value = -1880.2856
$ws.Range("H136").Value = 2578.75
$ws.Range("I136").Value = 2083.2856
$ws.Range("K136").Value = 6249.8568
$ws.Range("M136").Value = -3699.8568
$ws.Range("H141").Value = 60000
$ws.Range("J141").Value = 60000
$ws.Range("L141").Value = 60000
$ws.Range("N141").Value = -70360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1500
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 3000
$ws.Range("M68").Value = -2189
$ws.Range("H71").Value = 1500
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 9000
$ws.Range("M71").Value = -4944
$ws.Range("H107").Value = 271.5
$ws.Range("J107").Value = 271.5
$ws.Range("L107").Value = 814.5
$ws.Range("N107").Value = -4654.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 26257
$ws.Range("J43").Value = 26257
$ws.Range("L43").Value = 26257
$ws.Range("N43").Value = -26559
$ws.Range("H46").Value = 2234
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 3468
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 3468
$ws.Range("M46").Value = -844
$ws.Range("N46").Value = -3780
$ws.Range("H70").Value = 5274.4443
$ws.Range("I70").Value = 2750
$ws.Range("J70").Value = 5995.7144
$ws.Range("K70").Value = 2750
$ws.Range("L70").Value = 5995.7144
$ws.Range("M70").Value = -2480
$ws.Range("N70").Value = -6535.7144
$ws.Range("H73").Value = 5274.4443
$ws.Range("I73").Value = 2750
$ws.Range("J73").Value = 5995.7144
$ws.Range("K73").Value = 2750
$ws.Range("L73").Value = 5995.7144
$ws.Range("M73").Value = -1814
$ws.Range("N73").Value = -7867.7144
$ws.Range("H80").Value = 4664.6665
$ws.Range("I80").Value = 3596.4
$ws.Range("K80").Value = 3596.4
$ws.Range("M80").Value = -2598.4
$ws.Range("H83").Value = 4664.6665
$ws.Range("I83").Value = 3596.4
$ws.Range("K83").Value = 17982
$ws.Range("M83").Value = -12990
$ws.Range("H113").Value = 1313
$ws.Range("I113").Value = 1313
$ws.Range("K113").Value = 1313
$ws.Range("M113").Value = 857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2999.6667
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 2999.5
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 2999.5
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -3375.5
